$d = $word.ActiveDocument

# 1. Append " ERLEDIGT" to the paragraph ending with
#    "Die Anzahl der Gerichte verändern falls nötig"
$d.Content.Find.Execute(
    "Die Anzahl der Gerichte verändern falls nötig",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Die Anzahl der Gerichte verändern falls nötig ERLEDIGT",
    2
)

# 2. Replace " IMMER NOCH" with " ERLEDIGT"
#    (diff shows this split into two runs: " " and "ERLEDIGT",
#     but net visible text is " ERLEDIGT")
$d.Content.Find.Execute(
    " IMMER NOCH",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " ERLEDIGT",
    2
)
